$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.834.96"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.551.03"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.14%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "2.939.10"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.88%  "
$ws.Range("D16").Value = "2.617.83"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "42.837.10"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "0.0₃0955"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0805"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.112"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.03%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "1.988.72"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "2.791.70"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
